$d = $word.ActiveDocument

# Namespace declaration used for every InsertXML payload below.
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-RangeXml($range, $innerXml) {
    $payload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document ' + $wNs + '><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData>' +
        '</pkg:part></pkg:package>'
    $range.InsertXML($payload)
}

function Find-ParagraphContaining($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like $needle) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 0) The document carries a stray "_GoBack" bookmark left over from the last
#    editing session. Remove it now (before any other edit) so that later
#    InsertXML calls do not resurrect it in its old location.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 1) Clean up the now-bookmark-free empty paragraph that follows
#    "QA Testing Lead Carl Petersen".
# ---------------------------------------------------------------------------
$goBackHost = Find-ParagraphContaining("*Carl Petersen*")
$goBackHost = $goBackHost.Next()
Set-RangeXml $goBackHost.Range ('<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr></w:p>')

# ---------------------------------------------------------------------------
# 2) Shade the "Interpret feedback..." paragraph, re-home the "_GoBack"
#    bookmark at the very end of it, and replace the trailing empty
#    "ListParagraph" paragraph with the new task note paragraph. Both
#    paragraphs are rewritten together in a single InsertXML call because
#    the trailing paragraph is the last one in the body (its mark also
#    carries the section break) and replacing it alone would insert an
#    extra paragraph instead of overwriting it in place.
# ---------------------------------------------------------------------------
$interpretParagraph = Find-ParagraphContaining("*Interpret feedback from novice test-players*")
$taskParagraph = $interpretParagraph.Next()
$combinedRange = $d.Range($interpretParagraph.Range.Start, $taskParagraph.Range.End)

$taskText = "CARL NEEDS TO DO THIS: Test Report: " + [char]0x2022 + " Updated test plan (separate document) " + `
    [char]0x2022 + " Updates to any testing resources o Updates to automated testing processes o Updates to manual testing process documents " + `
    [char]0x2022 + " Summary of Issues being tracked o Number/severity of issues open o Number of new issues o Number of issues closed " + `
    [char]0x2022 + " Details of QA testing performed this sprint (when, where, who, what, how long) o Summary of testing results " + `
    [char]0x2022 + " Details of Play testing performed this sprint (when, where, who, what, how long) o Summary of testing results"

$taskTextEscaped = $taskText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")

$combinedInner =
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/></w:pPr><w:r><w:t xml:space="preserve">Interpret feedback from novice test-players </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' +
    '<w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/><w:ind w:left="1800"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="23"/><w:szCs w:val="23"/><w:shd w:val="clear" w:color="auto" w:fill="36393F"/></w:rPr><w:t>' + $taskTextEscaped + '</w:t></w:r></w:p>'

Set-RangeXml $combinedRange $combinedInner
